$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2881.25
$ws.Range("I32").Value = 1365.4
$ws.Range("J32").Value = 3964
$ws.Range("K32").Value = 1365.4
$ws.Range("L32").Value = 3964
$ws.Range("M32").Value = -1039.4
$ws.Range("N32").Value = -4616
$ws.Range("H46").Value = 2520
$ws.Range("J46").Value = 2520
$ws.Range("L46").Value = 7560
$ws.Range("N46").Value = -7798
$ws.Range("H60").Value = 2520
$ws.Range("J60").Value = 2520
$ws.Range("L60").Value = 7560
$ws.Range("N60").Value = -8528
$ws.Range("H64").Value = 4998.7144
$ws.Range("I64").Value = 4199.5
$ws.Range("K64").Value = 4199.5
$ws.Range("M64").Value = -3951.5
$ws.Range("H67").Value = 4998.7144
$ws.Range("I67").Value = 4199.5
$ws.Range("K67").Value = 4199.5
$ws.Range("M67").Value = -3341.5
$ws.Range("H75").Value = 45000
$ws.Range("J75").Value = 45000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46872
$ws.Range("H78").Value = 45000
$ws.Range("J78").Value = 45000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -144360
$ws.Range("H82").Value = 444.2
$ws.Range("I82").Value = 444.2
$ws.Range("K82").Value = 1332.6
$ws.Range("M82").Value = -926.5999999999999
$ws.Range("H85").Value = 444.2
$ws.Range("I85").Value = 444.2
$ws.Range("K85").Value = 1332.6
$ws.Range("M85").Value = 71.40000000000009
$ws.Range("H88").Value = 11459.2
$ws.Range("J88").Value = 11459.2
$ws.Range("L88").Value = 11459.2
$ws.Range("N88").Value = -12271.2
$ws.Range("H91").Value = 11459.2
$ws.Range("J91").Value = 11459.2
$ws.Range("L91").Value = 11459.2
$ws.Range("N91").Value = -14267.2
$ws.Range("H104").Value = 186.8
$ws.Range("I104").Value = 186.8
$ws.Range("K104").Value = 560.4000000000001
$ws.Range("M104").Value = 1186.6
$ws.Range("H107").Value = 269.875
$ws.Range("I107").Value = 291.7143
$ws.Range("J107").Value = 117
$ws.Range("K107").Value = 291.7143
$ws.Range("L107").Value = 117
$ws.Range("M107").Value = 1628.2857
$ws.Range("N107").Value = -3957
$ws.Range("H111").Value = 797.5
$ws.Range("J111").Value = 797.5
$ws.Range("L111").Value = 2392.5
$ws.Range("N111").Value = -8526.5
$ws.Range("H127").Value = 3434
$ws.Range("I127").Value = 4168
$ws.Range("K127").Value = 12504
$ws.Range("M127").Value = -7544

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 44933.332
$ws.Range("J43").Value = 44900
$ws.Range("L43").Value = 44900
$ws.Range("N43").Value = -45526
$ws.Range("H88").Value = 1553.5714
$ws.Range("I88").Value = 1497
$ws.Range("K88").Value = 1497
$ws.Range("M88").Value = -1091
$ws.Range("H91").Value = 1553.5714
$ws.Range("I91").Value = 1497
$ws.Range("K91").Value = 1497
$ws.Range("M91").Value = -93
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H124").Value = 66249.75
$ws.Range("J124").Value = 66249.75
$ws.Range("L124").Value = 66249.75
$ws.Range("N124").Value = -76069.75
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2132.9167
$ws.Range("I105").Value = 2288.4443
$ws.Range("J105").Value = 1666.3334
$ws.Range("K105").Value = 2288.4443
$ws.Range("L105").Value = 1666.3334
$ws.Range("M105").Value = -541.4443000000001
$ws.Range("N105").Value = -5160.3334
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H135").Value = 41235.758
$ws.Range("J135").Value = 41235.758
$ws.Range("L135").Value = 41235.758
$ws.Range("N135").Value = -51375.758

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 89097.38
$ws.Range("J94").Value = 4418
$ws.Range("L94").Value = 4418
$ws.Range("N94").Value = -5320

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 147.08333
$ws.Range("I2").Value = 56.5
$ws.Range("J2").Value = 192.375
$ws.Range("K2").Value = 339
$ws.Range("L2").Value = 1154.25
$ws.Range("M2").Value = -226
$ws.Range("N2").Value = -1380.25
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 296.44446
$ws.Range("J12").Value = 234.5
$ws.Range("L12").Value = 703.5
$ws.Range("N12").Value = -1049.5
$ws.Range("H109").Value = 922.7692
$ws.Range("I109").Value = 245.09091
$ws.Range("K109").Value = 735.27273
$ws.Range("M109").Value = 304.72727
$ws.Range("H114").Value = 2038.5555
$ws.Range("I114").Value = 550
$ws.Range("J114").Value = 2463.8572
$ws.Range("K114").Value = 1650
$ws.Range("L114").Value = 7391.571599999999
$ws.Range("M114").Value = 1604
$ws.Range("N114").Value = -13899.5716
$ws.Range("H121").Value = 9305.223
$ws.Range("I121").Value = 27642.5
$ws.Range("J121").Value = 4066
$ws.Range("K121").Value = 82927.5
$ws.Range("L121").Value = 12198
$ws.Range("M121").Value = -81617.5
$ws.Range("N121").Value = -14818

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 49704.332
$ws.Range("J15").Value = 49704.332
$ws.Range("L15").Value = 49704.332
$ws.Range("N15").Value = -50280.332
$ws.Range("H81").Value = 49704.332
$ws.Range("J81").Value = 49704.332
$ws.Range("L81").Value = 49704.332
$ws.Range("N81").Value = -51700.332
$ws.Range("H84").Value = 49704.332
$ws.Range("J84").Value = 49704.332
$ws.Range("L84").Value = 149112.996
$ws.Range("N84").Value = -159096.996
$ws.Range("H132").Value = 3655.75
$ws.Range("I132").Value = 3655.75
$ws.Range("K132").Value = 10967.25
$ws.Range("M132").Value = -8437.25

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1912.2727
$ws.Range("I22").Value = 1424.7
$ws.Range("K22").Value = 1424.7
$ws.Range("M22").Value = -1129.7
$ws.Range("H27").Value = 1912.2727
$ws.Range("I27").Value = 1424.7
$ws.Range("K27").Value = 1424.7
$ws.Range("M27").Value = -1317.7
$ws.Range("H82").Value = 2591.6667
$ws.Range("I82").Value = 2816.6667
$ws.Range("J82").Value = 2141.6667
$ws.Range("K82").Value = 2816.6667
$ws.Range("L82").Value = 2141.6667
$ws.Range("M82").Value = -2455.6667
$ws.Range("N82").Value = -2863.6667
$ws.Range("H85").Value = 2591.6667
$ws.Range("I85").Value = 2816.6667
$ws.Range("J85").Value = 2141.6667
$ws.Range("K85").Value = 2816.6667
$ws.Range("L85").Value = 2141.6667
$ws.Range("M85").Value = -1568.6667
$ws.Range("N85").Value = -4637.6667
$ws.Range("H132").Value = 9860.125
$ws.Range("I132").Value = 10468.714
$ws.Range("K132").Value = 31406.142
$ws.Range("M132").Value = -28876.142

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 17
$ws.Range("I31").Value = 17
$ws.Range("K31").Value = 17
$ws.Range("M31").Value = 331
$ws.Range("H81").Value = 1002118.1
$ws.Range("I81").Value = 1435
$ws.Range("J81").Value = 2002801.2
$ws.Range("K81").Value = 2870
$ws.Range("L81").Value = 4005602.4
$ws.Range("M81").Value = -1809
$ws.Range("N81").Value = -4007724.4
$ws.Range("H84").Value = 1002118.1
$ws.Range("I84").Value = 1435
$ws.Range("J84").Value = 2002801.2
$ws.Range("K84").Value = 14350
$ws.Range("L84").Value = 20028012
$ws.Range("M84").Value = -9046
$ws.Range("N84").Value = -20038620
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
